# TC_222 - Updated test data for normal load, cable capacitance etc.
# Adds "Instrinsically Safe Unit Details" section header (E1/E2) to the
# "Add Devices Loop A" sheet and moves the active selection to E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E1: header-style cell (copy formatting from A2, which already uses the
# bold/bordered header style), then set its text.
$ws.Range("A2").Copy()
$ws.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E1").Value = "Instrinsically Safe Unit Details"

# E2: shaded data-style cell (copy formatting from C4, which uses the grey
# fill + border style), then set its text.
$ws.Range("C4").Copy()
$ws.Range("E2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E2").Value = "Built-in Loop-A Intrinsically-safe Units"

$excel.CutCopyMode = $false

# Move/record the active selection as it was left in the saved workbook.
$ws.Range("E4").Select()
